$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Block 1: rows 14-18 ----
$ws.Range("H14").Value = 0.25

$ws.Range("A15").Value = "MIN"
$ws.Range("B15").Value = 150
$ws.Range("C15").Value = 160
$ws.Range("D15").Value = 0.05
$ws.Range("E15").Formula = "=C15*(1-D15)"
$ws.Range("G15").Value = 0.05
$ws.Range("H15").Formula = "=C15*(1-G15+`$H`$14)"

$ws.Range("A16").Value = "MAX"
$ws.Range("B16").Value = 231
$ws.Range("C16").Value = 275
$ws.Range("D16").Value = 0.4
$ws.Range("G16").Value = 0.5

$ws.Range("C17").Value = 325
$ws.Range("D17").Value = 0.45
$ws.Range("G17").Value = 0.55000000000000004

$ws.Range("C18").Value = 425
$ws.Range("D18").Value = 0.5
$ws.Range("G18").Value = 0.6

$ws.Range("E16:E18").Formula = "=C16*(1-D16)"
$ws.Range("H16:H18").Formula = "=C16*(1-G16+`$H`$14)"

# ---- Block 2: rows 21-25 ----
$ws.Range("H21").Value = 0.25

$ws.Range("A22").Value = "MIN"
$ws.Range("B22").Value = 150
$ws.Range("C22").Value = 160
$ws.Range("D22").Value = 0.1
$ws.Range("E22").Formula = "=C22*(1-D22)"
$ws.Range("G22").Value = 0.3
$ws.Range("H22").Formula = "=C22*(1-G22+`$H`$14)"

$ws.Range("A23").Value = "MAX"
$ws.Range("B23").Value = 154
$ws.Range("C23").Value = 260
$ws.Range("D23").Value = 0.42
$ws.Range("G23").Value = 0.66

$ws.Range("C24").Value = 310
$ws.Range("D24").Value = 0.51
$ws.Range("G24").Value = 0.755

$ws.Range("C25").Value = 410
$ws.Range("D25").Value = 0.6
$ws.Range("G25").Value = 0.8

$ws.Range("E23:E25").Formula = "=C23*(1-D23)"
$ws.Range("H23:H25").Formula = "=C23*(1-G23+`$H`$14)"

# ---- View changes ----
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("J21").Select()
